$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (2023-09-09 -> 2023-09-10, i.e. serial 45178 -> 45179) for every data row
# (rows 2 through 260).
$ws.Range("C2:C260").Value = 45179
